# "Add files via upload" — populate the subtitle placeholder on the title
# slide with the homework archive name and a hyperlinked Google Drive
# folder link, followed by a trailing empty paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder (ppPlaceholderSubtitle = 4) rather than
# assuming a fixed shape index.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.PlaceholderFormat.Type -eq 4) {
        $sh = $candidate
        break
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(2)
}

$tr = $sh.TextFrame.TextRange

# Build the text run-by-run so each insertion becomes its own <a:r>,
# mirroring how PowerPoint records incremental typed edits.
$tr.Text = "h"
[void]$tr.InsertAfter("w01-04")
[void]$tr.InsertAfter(".ppt")
[void]$tr.InsertAfter("：")

$linkText = "https://drive.google.com/drive/folders/1-pDV9bIA_EX6QU-Ci_Drk7dGkvZ8RmOG?usp=drive_link"
[void]$tr.InsertAfter($linkText)

# A trailing, otherwise-empty second paragraph follows the link line.
[void]$tr.InsertAfter("`r")

# Turn just the URL run into a hyperlink pointing at the same address
# (found by locating it in the now-complete text, so the paragraph break
# above doesn't inherit the link formatting).
$full = $sh.TextFrame.TextRange
$fullText = $full.Text
$linkStart = $fullText.IndexOf("https://") + 1
$linkRange = $full.Characters($linkStart, $linkText.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = $linkText
